$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D" (Price) column holds numeric-looking text such as "1.001" or
# "20.099.83". Force the cell number format to Text ("@") before writing
# so Excel keeps the exact original string instead of coercing it to a
# number (which would also collapse formats like "1.000" -> 1).

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.45"
$ws.Range("E35").Value = "  +7.56%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.367"
$ws.Range("E36").Value = "  -6.72%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05705"
$ws.Range("E38").Value = "  -6.91%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.730"
$ws.Range("E39").Value = "  -7.72%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.099.83"
$ws.Range("E2").Value = "  -7.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.428.71"
$ws.Range("E3").Value = "  -7.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.93"
$ws.Range("E6").Value = "  -5.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3730"
$ws.Range("E7").Value = "  -4.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3093"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.22"
$ws.Range("E9").Value = "  -6.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06594"
$ws.Range("E11").Value = "  -8.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.417"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.28"
$ws.Range("E14").Value = "  -7.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.181"
$ws.Range("E15").Value = "  -6.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.430.08"
$ws.Range("E16").Value = "  -7.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001013"
$ws.Range("E17").Value = "  -8.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05825"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.71"
$ws.Range("E19").Value = "  -8.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.687"
$ws.Range("E21").Value = "  -7.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.50"
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.333"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.118.46"
$ws.Range("E25").Value = "  -7.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.290"
$ws.Range("E26").Value = "  -4.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "138.37"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.95"
$ws.Range("E28").Value = "  -7.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.591.31"
$ws.Range("E29").Value = "  -7.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.53"
$ws.Range("E30").Value = "  -7.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.938"
$ws.Range("E31").Value = "  -18.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9093"
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.439"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07782"
$ws.Range("E34").Value = "  -5.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1923"
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.124"
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02016"
$ws.Range("E42").Value = "  -8.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.303"
$ws.Range("E43").Value = "  -8.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5349"
$ws.Range("E44").Value = "  -7.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.550"
$ws.Range("E45").Value = "  -5.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  -6.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5155"
$ws.Range("E47").Value = "  -6.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.781"
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.68"
$ws.Range("E49").Value = "  -5.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.054"
$ws.Range("E50").Value = "  -6.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  +0.00%  "
